# Finalized expense tracker project
# Rename header "Source" -> "Category", update the existing expense row,
# and append the rest of the expense entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Category"

# --- Existing row (2) gets overwritten with the "Rent" expense ---
$ws.Range("A2").Value = "Rent"
$ws.Range("B2").Value = 800
$ws.Range("C2").Value = 45930.33362268518

# --- New expense rows (3-6) ---
$data = @(
    @("Groceries", 250, 45905.33362268518),
    @("Travel", 1000, 45903.33362268518),
    @("Loan", 400, 45903.33362268518),
    @("Spotify Subscription", 150, 45902.33362268518)
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# The date column (C) needs the same date number format the original C2
# cell already carried (style index 1, numFmtId 14) - copy formats only so
# no new style/numFmt entries get minted.
$ws.Range("C2").Copy()
$ws.Range("C3:C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
